$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.115.13"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.122.01"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.36"
$ws.Range("E5").Value = "  -2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.66"
$ws.Range("E6").Value = "  -5.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.118.97"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("E11").Value = "  -3.96%  "
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.19"
$ws.Range("E14").Value = "  -3.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.634.52"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.018.33"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.120.23"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.79"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.69"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.28"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.05"
$ws.Range("E25").Value = "  -4.78%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").Value = "  -7.41%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.95"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("E33").Value = "  -8.06%  "
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.99"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0716"
$ws.Range("E38").Value = "  -4.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0389"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.42"
$ws.Range("E40").Value = "  -7.50%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("E43").Value = "  -11.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.884.65"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  -5.35%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.85"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  -6.85%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.56"
$ws.Range("E51").Value = "  -1.99%  "
